# This blood-test report sheet was populated from noisy OCR text, where the
# "result / reference range / unit" columns were smashed together and row
# labels carried a stray leading index digit. This script cleans each row:
#   - strips the leading index number baked into column A's label,
#   - splits the jammed "range+unit" text in column C into a clean
#     "low-high" range, and
#   - puts the correct unit into column D.
# A handful of result values in column B are also corrected/re-typed;
# because several of those look like plain numbers ("122.0", "307", "0.1",
# ...), a leading apostrophe is used so Excel keeps them as text (matching
# the original text-formatted lab values) instead of silently converting
# them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 白细胞 (White blood cell count)
$ws.Range('A2').Value = '白细胞'
$ws.Range('C2').Value = '4.00-10.0'
$ws.Range('D2').Value = '10^9/L'

# Row 3 - 红细胞 (Red blood cell count)
$ws.Range('A3').Value = '红细胞'
$ws.Range('C3').Value = '3.50-5.50'
$ws.Range('D3').Value = '10^12/L'

# Row 4 - 血红蛋白 (Hemoglobin)
$ws.Range('A4').Value = '血红蛋白'
$ws.Range('B4').Value = "'122.0"
$ws.Range('C4').Value = '010.0-160.'
$ws.Range('D4').Value = 'g/L'

# Row 5 - 红细胞压积 (Hematocrit)
$ws.Range('A5').Value = '红细胞压积'
$ws.Range('B5').Value = "'35.0"
$ws.Range('C5').Value = '436.0-50.0'
$ws.Range('D5').Value = '%'

# Row 6 - 血小板 (Platelet count); result value cleared entirely
$ws.Range('A6').Value = '血小板'
$ws.Range('B6').ClearContents()
$ws.Range('C6').Value = "'307"
$ws.Range('D6').Value = '10^9/L'

# Row 7 - (MPV) reference range / unit only
$ws.Range('C7').Value = '9.0-13.0'
$ws.Range('D7').Value = 'fL'

# Row 8 - 血小板压积 (Plateletcrit)
$ws.Range('A8').Value = '血小板压积'

# Row 9 - 红细胞平均体积 (Mean corpuscular volume)
$ws.Range('A9').Value = '红细胞平均体积'
$ws.Range('B9').Value = "'79.0"
$ws.Range('C9').Value = '6486.0-100.'
$ws.Range('D9').Value = 'fL'

# Row 10 - 平均血红蛋白量 (Mean corpuscular hemoglobin)
$ws.Range('A10').Value = '平均血红蛋白量'
$ws.Range('C10').Value = '26.0-33.0'
$ws.Range('D10').Value = 'pg'

# Row 11 - 平均血红蛋白浓度 (Mean corpuscular hemoglobin concentration)
$ws.Range('A11').Value = '平均血红蛋白浓度'
$ws.Range('C11').Value = '310-370'
$ws.Range('D11').Value = 'g/L'

# Row 12 - 中性细胞比率 (Neutrophil ratio)
$ws.Range('A12').Value = '中性细胞比率'
$ws.Range('B12').Value = "'38.0"
$ws.Range('C12').Value = '445.0-77.0'
$ws.Range('D12').Value = '%'

# Row 13 - 淋巴细胞比率 (Lymphocyte ratio)
$ws.Range('A13').Value = '淋巴细胞比率'
$ws.Range('B13').Value = "'50.4"
$ws.Range('C13').Value = '20.0-40.0'
$ws.Range('D13').Value = '%'

# Row 14 - 单核细胞比率 (Monocyte ratio)
$ws.Range('A14').Value = '单核细胞比率'
$ws.Range('C14').Value = '3.0-8.0'
$ws.Range('D14').Value = '%'

# Row 15 - 嗜酸性粒细胞比率 (Eosinophil ratio)
$ws.Range('A15').Value = '嗜酸性粒细胞比率'
$ws.Range('C15').Value = '0.5-5.0'
$ws.Range('D15').Value = '%'

# Row 16 - 嗜碱性粒细胞比率 (Basophil ratio)
$ws.Range('A16').Value = '嗜碱性粒细胞比率'
$ws.Range('C16').Value = '0.0-1.0'

# Row 17 - 中性细胞数 (Neutrophil count)
$ws.Range('A17').Value = '中性细胞数'
$ws.Range('C17').Value = '2.0-7.7109'

# Row 18 - 淋巴细胞数 (Lymphocyte count)
$ws.Range('A18').Value = '淋巴细胞数'
$ws.Range('C18').Value = '0.80-4.00109'

# Row 19 - 单核细胞数 (Monocyte count)
$ws.Range('A19').Value = '单核细胞数'
$ws.Range('C19').Value = '0.12-0.80109'

# Row 20 - 嗜酸性粒细胞数 (Eosinophil count)
$ws.Range('A20').Value = '嗜酸性粒细胞数'
$ws.Range('B20').Value = "'0.1"
$ws.Range('C20').Value = '0.05-0.50109'

# Row 21 - 嗜碱性粒细胞 (Basophil count)
$ws.Range('A21').Value = '嗜碱性粒细胞'
$ws.Range('B21').Value = "'0.0"
$ws.Range('C21').Value = '0.00-0.10109'

# Row 22 - (HCT-related) reference range only
$ws.Range('C22').Value = '37.0-50.0'

# Row 23 - RDW-CV
$ws.Range('A23').Value = 'RDW-CV'

# Row 24 - 血小板分布宽度 (Platelet distribution width)
$ws.Range('A24').Value = '血小板分布宽度'
$ws.Range('C24').Value = '9.0-17.0'

# Row 25 - 大型血小板比率 (Platelet-large cell ratio)
$ws.Range('A25').Value = '大型血小板比率'
$ws.Range('C25').Value = '13.0-43.0'
